$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.830.19'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '3.194.12'
$ws.Range('E3').Value = '  -3.07%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'213.78"
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = "'618.11"
$ws.Range('E6').Value = '  -1.77%  '
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('D8').Value = "'0.691"
$ws.Range('E8').Value = '  -3.81%  '
$ws.Range('D9').Value = "'1.00"
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = '3.189.10'
$ws.Range('E10').Value = '  -3.08%  '
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('E12').Value = '  -5.32%  '
$ws.Range('E13').Value = '  -4.25%  '
$ws.Range('D14').Value = '90.584.57'
$ws.Range('E14').Value = '  +1.44%  '
$ws.Range('D15').Value = '3.771.59'
$ws.Range('E15').Value = '  -3.02%  '
$ws.Range('D16').Value = "'32.90"
$ws.Range('E16').Value = '  -4.00%  '
$ws.Range('D17').Value = "'5.25"
$ws.Range('E17').Value = '  -3.75%  '
$ws.Range('D18').Value = '3.190.33'
$ws.Range('E18').Value = '  -3.41%  '
$ws.Range('D19').Value = "'3.27"
$ws.Range('E19').Value = '  +6.18%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'446.70"
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = "'13.51"
$ws.Range('E21').Value = '  -4.43%  '
$ws.Range('D22').Value = "'0.0000188"
$ws.Range('E22').Value = '  +35.71%  '
$ws.Range('D23').Value = "'8.60"
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('E24').Value = '  -4.62%  '
$ws.Range('D25').Value = "'5.17"
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('D26').Value = "'11.81"
$ws.Range('E26').Value = '  -3.36%  '
$ws.Range('D27').Value = '3.353.70'
$ws.Range('E27').Value = '  -3.20%  '
$ws.Range('D28').Value = "'74.98"
$ws.Range('E28').Value = '  -2.74%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  -7.19%  '
$ws.Range('D31').Value = "'0.998"
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('D32').Value = "'4.24"
$ws.Range('E32').Value = '  +36.33%  '
$ws.Range('D33').Value = "'8.47"
$ws.Range('E33').Value = '  -4.73%  '
$ws.Range('D34').Value = "'535.01"
$ws.Range('E34').Value = '  -5.16%  '
$ws.Range('D35').Value = "'7.01"
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('E36').Value = '  -4.89%  '
$ws.Range('D37').Value = "'1.26"
$ws.Range('E37').Value = '  -8.75%  '
$ws.Range('D38').Value = "'22.00"
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('D39').Value = "'22.33"
$ws.Range('E39').Value = '  +2.23%  '
$ws.Range('E40').Value = '  -8.17%  '
$ws.Range('D41').Value = "'0.999"
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D43').Value = "'0.377"
$ws.Range('E43').Value = '  -6.51%  '
$ws.Range('D44').Value = "'1.92"
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('D45').Value = "'147.70"
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('D46').Value = "'44.39"
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('D47').Value = "'172.71"
$ws.Range('E47').Value = '  -4.76%  '
$ws.Range('D48').Value = "'0.124"
$ws.Range('E48').Value = '  -7.06%  '
$ws.Range('D49').Value = "'1.25"
$ws.Range('E49').Value = '  -5.37%  '
$ws.Range('D50').Value = "'0.614"
$ws.Range('E50').Value = '  -2.27%  '
$ws.Range('E51').Value = '  -4.00%  '
